$d = $word.ActiveDocument

# The page footer boilerplate ("Ver no Jupiter ..." and the copyright line),
# plus the blank separator paragraph that precedes them, were removed from
# the rebuilt site page. Locate those two paragraphs by their text and
# delete the whole run, including the blank paragraph right before them.

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$count = $d.Paragraphs.Count
$idx1 = -1
$idx2 = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq $target1) { $idx1 = $i }
    if ($t -eq $target2) { $idx2 = $i }
}

if ($idx1 -gt 0 -and $idx2 -gt $idx1) {
    # Include the blank paragraph immediately preceding the "Ver no Jupiter" line.
    $startIdx = $idx1 - 1
    $pStart = $d.Paragraphs.Item($startIdx)
    # End right at the start of the paragraph following the copyright line so
    # that paragraph's own properties are left completely untouched.
    $pAfter = $d.Paragraphs.Item($idx2 + 1)

    $rng = $d.Range($pStart.Range.Start, $pAfter.Range.Start)
    $rng.Delete()
}
